$d = $word.ActiveDocument

# Each three-digit-divided-by-one-digit expression in the document is unique,
# so we can safely use Find & Replace (exact match) for each pair, one at a time,
# re-searching the whole document content range for every replacement.
$pairs = @(
    ,@("476÷2=", "598÷3=")
    ,@("309÷7=", "612÷4=")
    ,@("359÷3=", "306÷3=")
    ,@("844÷4=", "770÷8=")
    ,@("699÷5=", "761÷4=")
    ,@("894÷4=", "316÷9=")
    ,@("241÷4=", "118÷7=")
    ,@("124÷4=", "536÷7=")
    ,@("181÷7=", "886÷2=")
    ,@("397÷2=", "452÷3=")
    ,@("234÷4=", "920÷3=")
    ,@("196÷4=", "179÷9=")
    ,@("100÷6=", "215÷8=")
    ,@("302÷2=", "162÷6=")
    ,@("791÷2=", "764÷8=")
    ,@("103÷9=", "247÷4=")
    ,@("802÷9=", "317÷3=")
    ,@("480÷7=", "579÷5=")
    ,@("629÷8=", "516÷6=")
    ,@("957÷2=", "724÷9=")
    ,@("678÷3=", "496÷7=")
    ,@("615÷7=", "621÷5=")
    ,@("328÷5=", "914÷7=")
    ,@("272÷5=", "823÷2=")
    ,@("390÷2=", "210÷3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: could not find $old"
    }
}

Write-Output "Done replacing $($pairs.Count) expressions."
